{"js": "// Replace each arithmetic-equation cell in the (sole) table with its new\n// value, in row-major order, matching the OOXML diff exactly. Positional\n// (row,col) addressing is used instead of text search-and-replace because\n// several \"before\" equations are duplicated (e.g. \"49+18=67\" appears twice)\n// but map to different \"after\" values depending on position.\nconst NEW_VALUES = [\n  [\"86-45=41\", \"54+27=81\", \"28-18=10\", \"72-31=41\", \"22+28=50\"],\n  [\"86-82=4\", \"34+15=49\", \"17+27=44\", \"10-1=9\", \"8+32=40\"],\n  [\"28+48=76\", \"89-80=9\", \"65+25=90\", \"16+83=99\", \"15+80=95\"],\n  [\"78-56=22\", \"89-39=50\", \"34-25=9\", \"53-45=8\", \"22+76=98\"],\n  [\"34-11=23\", \"4+43=47\", \"81-24=57\", \"29+53=82\", \"86+11=97\"],\n  [\"91-7=84\", \"19-4=15\", \"65-32=33\", \"87-13=74\", \"62-42=20\"],\n  [\"19-10=9\", \"39+13=52\", \"51-29=22\", \"36-12=24\", \"60+8=68\"],\n  [\"50-29=21\", \"75-25=50\", \"16+82=98\", \"95-24=71\", \"87-56=31\"],\n  [\"94-69=25\", \"95+2=97\", \"22+4=26\", \"62+35=97\", \"94-75=19\"],\n  [\"41-18=23\", \"89-15=74\", \"80-74=6\", \"93-81=12\", \"0+25=25\"],\n  [\"15+23=38\", \"34+50=84\", \"19+38=57\", \"80+7=87\", \"9+60=69\"],\n  [\"22+60=82\", \"47+15=62\", \"24+17=41\", \"51+2=53\", \"62-40=22\"],\n  [\"39+35=74\", \"21+21=42\", \"50+38=88\", \"28-18=10\", \"94-4=90\"],\n  [\"85+5=90\", \"67+20=87\", \"23+74=97\", \"82-52=30\", \"21+75=96\"],\n  [\"26-3=23\", \"97-20=77\", \"41-33=8\", \"4+48=52\", \"0+33=33\"],\n  [\"33+26=59\", \"45+36=81\", \"17-1=16\", \"3+19=22\", \"23+58=81\"],\n  [\"77-55=22\", \"87-62=25\", \"94-69=25\", \"22-6=16\", \"48+30=78\"],\n  [\"90+5=95\", \"79+9=88\", \"42+13=55\", \"96-33=63\", \"30+63=93\"],\n  [\"21+38=59\", \"4+87=91\", \"63+32=95\", \"31-0=31\", \"35+7=42\"],\n  [\"58-22=36\", \"81-22=59\", \"32+62=94\", \"46-39=7\", \"90-67=23\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.rowCount !== NEW_VALUES.length) {\n  throw new Error(\n    `Expected ${NEW_VALUES.length} rows, found ${table.rowCount}`\n  );\n}\n\nfor (let r = 0; r < NEW_VALUES.length; r++) {\n  for (let c = 0; c < NEW_VALUES[r].length; c++) {\n    table.getCell(r, c).value = NEW_VALUES[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-equation cell in the (sole) table with its new\n# value, in row-major order, matching the OOXML diff exactly. Positional\n# (row,col) addressing is used instead of text search-and-replace because\n# several \"before\" equations are duplicated (e.g. \"49+18=67\" appears twice)\n# but map to different \"after\" values depending on position.\n$newValues = @(\n    @(\"86-45=41\", \"54+27=81\", \"28-18=10\", \"72-31=41\", \"22+28=50\"),\n    @(\"86-82=4\", \"34+15=49\", \"17+27=44\", \"10-1=9\", \"8+32=40\"),\n    @(\"28+48=76\", \"89-80=9\", \"65+25=90\", \"16+83=99\", \"15+80=95\"),\n    @(\"78-56=22\", \"89-39=50\", \"34-25=9\", \"53-45=8\", \"22+76=98\"),\n    @(\"34-11=23\", \"4+43=47\", \"81-24=57\", \"29+53=82\", \"86+11=97\"),\n    @(\"91-7=84\", \"19-4=15\", \"65-32=33\", \"87-13=74\", \"62-42=20\"),\n    @(\"19-10=9\", \"39+13=52\", \"51-29=22\", \"36-12=24\", \"60+8=68\"),\n    @(\"50-29=21\", \"75-25=50\", \"16+82=98\", \"95-24=71\", \"87-56=31\"),\n    @(\"94-69=25\", \"95+2=97\", \"22+4=26\", \"62+35=97\", \"94-75=19\"),\n    @(\"41-18=23\", \"89-15=74\", \"80-74=6\", \"93-81=12\", \"0+25=25\"),\n    @(\"15+23=38\", \"34+50=84\", \"19+38=57\", \"80+7=87\", \"9+60=69\"),\n    @(\"22+60=82\", \"47+15=62\", \"24+17=41\", \"51+2=53\", \"62-40=22\"),\n    @(\"39+35=74\", \"21+21=42\", \"50+38=88\", \"28-18=10\", \"94-4=90\"),\n    @(\"85+5=90\", \"67+20=87\", \"23+74=97\", \"82-52=30\", \"21+75=96\"),\n    @(\"26-3=23\", \"97-20=77\", \"41-33=8\", \"4+48=52\", \"0+33=33\"),\n    @(\"33+26=59\", \"45+36=81\", \"17-1=16\", \"3+19=22\", \"23+58=81\"),\n    @(\"77-55=22\", \"87-62=25\", \"94-69=25\", \"22-6=16\", \"48+30=78\"),\n    @(\"90+5=95\", \"79+9=88\", \"42+13=55\", \"96-33=63\", \"30+63=93\"),\n    @(\"21+38=59\", \"4+87=91\", \"63+32=95\", \"31-0=31\", \"35+7=42\"),\n    @(\"58-22=36\", \"81-22=59\", \"32+62=94\", \"46-39=7\", \"90-67=23\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nif ($t.Rows.Count -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) rows, found $($t.Rows.Count)\"\n}\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
